$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D and E to be treated as literal text while we assign the
# new values, so that strings such as "30.459.12" or "1.210" are not
# auto-coerced into numbers (which would lose formatting/precision).
$dRange = $ws.Range("D2:D51")
$eRange = $ws.Range("E2:E51")
$dRange.NumberFormat = "@"
$eRange.NumberFormat = "@"

$ws.Range("D2").Value = "30.459.12"
$ws.Range("E2").Value = "  -0.72%  "

$ws.Range("D3").Value = "2.105.67"
$ws.Range("E3").Value = "  +0.14%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").Value = "336.37"
$ws.Range("E5").Value = "  +1.92%  "

$ws.Range("E6").Value = "  +0.06%  "

$ws.Range("E7").Value = "  -0.47%  "

$ws.Range("D8").Value = "0.4612"
$ws.Range("E8").Value = "  +6.49%  "

$ws.Range("D9").Value = "52.33"
$ws.Range("E9").Value = "  +12.21%  "

$ws.Range("D10").Value = "0.08952"
$ws.Range("E10").Value = "  +0.32%  "

$ws.Range("D11").Value = "1.179"
$ws.Range("E11").Value = "  +1.04%  "

$ws.Range("D12").Value = "24.46"
$ws.Range("E12").Value = "  -0.56%  "

$ws.Range("D13").Value = "2.092.12"
$ws.Range("E13").Value = "  -0.44%  "

$ws.Range("D14").Value = "6.796"
$ws.Range("E14").Value = "  +1.58%  "

$ws.Range("D15").Value = "7.892"
$ws.Range("E15").Value = "  +1.30%  "

$ws.Range("D16").Value = "96.39"
$ws.Range("E16").Value = "  -0.50%  "

$ws.Range("E17").Value = "  +0.06%  "

$ws.Range("D18").Value = "0.00001131"
$ws.Range("E18").Value = "  +0.52%  "

$ws.Range("D19").Value = "0.06627"
$ws.Range("E19").Value = "  -0.44%  "

$ws.Range("D20").Value = "19.30"
$ws.Range("E20").Value = "  +1.96%  "

$ws.Range("E21").Value = "  +0.00%  "

$ws.Range("D22").Value = "6.284"
$ws.Range("E22").Value = "  -0.13%  "

$ws.Range("D23").Value = "30.522.82"
$ws.Range("E23").Value = "  -0.72%  "

$ws.Range("D24").Value = "12.34"
$ws.Range("E24").Value = "  +0.35%  "

$ws.Range("E25").Value = "  +3.48%  "

$ws.Range("D26").Value = "2.340.09"
$ws.Range("E26").Value = "  -0.36%  "

$ws.Range("D27").Value = "22.31"
$ws.Range("E27").Value = "  -0.70%  "

$ws.Range("D28").Value = "2.571"
$ws.Range("E28").Value = "  +0.30%  "

$ws.Range("D29").Value = "163.58"
$ws.Range("E29").Value = "  +1.08%  "

$ws.Range("D30").Value = "132.72"
$ws.Range("E30").Value = "  -0.06%  "

$ws.Range("D31").Value = "1.199"
$ws.Range("E31").Value = "  -0.47%  "

$ws.Range("E32").Value = "  -0.20%  "

$ws.Range("E33").Value = "  +9.63%  "

$ws.Range("D34").Value = "6.162"
$ws.Range("E34").Value = "  +0.31%  "

$ws.Range("D35").Value = "3.922"
$ws.Range("E35").Value = "  +1.76%  "

$ws.Range("D36").Value = "10.44"
$ws.Range("E36").Value = "  +7.94%  "

$ws.Range("D37").Value = "0.02571"
$ws.Range("E37").Value = "  -0.69%  "

$ws.Range("D38").Value = "0.06834"
$ws.Range("E38").Value = "  +1.74%  "

$ws.Range("D39").Value = "5.544"
$ws.Range("E39").Value = "  +0.35%  "

$ws.Range("D40").Value = "12.85"
$ws.Range("E40").Value = "  +1.95%  "

$ws.Range("E41").Value = "  +0.78%  "

$ws.Range("D42").Value = "0.6898"
$ws.Range("E42").Value = "  +1.56%  "

$ws.Range("D43").Value = "1.245"
$ws.Range("E43").Value = "  -0.31%  "

$ws.Range("E44").Value = "  +5.65%  "

$ws.Range("E45").Value = "  -0.02%  "

$ws.Range("D46").Value = "0.6389"
$ws.Range("E46").Value = "  +0.10%  "

$ws.Range("E47").Value = "  -0.20%  "

$ws.Range("D48").Value = "0.00000000363"
$ws.Range("E48").Value = "  +25.19%  "

$ws.Range("D49").Value = "3.666"
$ws.Range("E49").Value = "  +0.97%  "

$ws.Range("E50").Value = "  -0.37%  "

$ws.Range("D51").Value = "1.210"
$ws.Range("E51").Value = "  +1.35%  "

# Restore the default (General) number format / style so the cells keep
# looking like the rest of the sheet (no visible formatting change).
$dRange.NumberFormat = "General"
$eRange.NumberFormat = "General"
$dRange.Style = "Normal"
$eRange.Style = "Normal"
